$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    3  = @("三花智控","三花智控","龙洲股份")
    4  = @("海王生物","和而泰","海王生物")
    5  = @("和而泰","海王生物","平潭发展")
    6  = @("平潭发展","平潭发展","和而泰")
    7  = @("航天机电","巨轮智能","实达集团")
    8  = @("中国卫星","航天机电","三花智控")
    9  = @("巨轮智能","中国卫星","航天机电")
    10 = @("合富中国","合富中国","道明光学")
    11 = @("海欣食品","华映科技","合富中国")
    12 = @("龙洲股份","大众公用","中国铀业")
    13 = @("C中国铀","海欣食品","大众公用")
    14 = @("睿能科技","C中国铀","海欣食品")
    15 = @("大众公用","龙洲股份","巨轮智能")
    16 = @("实达集团","达华智能","太阳电缆")
    17 = @("达华智能","华伍股份","达华智能")
    18 = @("福蓉科技","航天电子","国机重装")
    19 = @("龙溪股份","日发精机","雷科防务")
    20 = @("航天动力","福蓉科技","四川金顶")
    21 = @("雷科防务","龙溪股份","中国卫星")
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
}
